$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data of row 2 and row 3 (the two observation records
# change places, e.g. after a re-sort/re-import), while every other column
# (C, D, J, K, N, P, T, U, V, W, Y, Z, AA, AB, AD, AE, AF, AG, AT, AW, AX, AY)
# holds identical values in both rows, so only A, B, E, F, G, H, I, Q, R, S
# (plus the stray empty L cell) actually need to move.

# --- Row 2 becomes what row 3 used to hold ---
$ws.Range("A2").Value = 111799186
$ws.Range("B2").Value = 89183
$ws.Range("E2").Value = 3215
$ws.Range("F2").Value = "Rödgul trumpetsvamp"
$ws.Range("G2").Value = "Craterellus lutescens"
$ws.Range("H2").Value = "(Fr.) Fr."

# Column I is stored as text ("20"/"1"), not a number - force text formatting
# before writing so the literal is kept as a string, then drop back to the
# default style so no stray formatting is left behind.
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "20"
$ws.Range("I2").Style = "Normal"

$ws.Range("Q2").Value = 513784.975650124
$ws.Range("R2").Value = 6704706.693730025
$ws.Range("S2").Value = 10

# L2 held an (empty) cell in the original row 2 but row 3 has none there -
# after the swap, row 2 should no longer have a populated L cell.
$ws.Range("L2").ClearContents()

# --- Row 3 becomes what row 2 used to hold ---
$ws.Range("A3").Value = 111799311
$ws.Range("B3").Value = 96251
$ws.Range("E3").Value = 220093
$ws.Range("F3").Value = "Korallrot"
$ws.Range("G3").Value = "Corallorhiza trifida"
$ws.Range("H3").Value = "Châtel."

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "1"
$ws.Range("I3").Style = "Normal"

$ws.Range("Q3").Value = 513717.9300395954
$ws.Range("R3").Value = 6704676.858456986
$ws.Range("S3").Value = 25

# Row 3 gains the (empty) L cell that row 2 used to carry.
$ws.Range("L3").NumberFormat = "@"
$ws.Range("L3").Value = ""
$ws.Range("L3").Style = "Normal"
